# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AB1) onto the
# new header cells so they match the look of the other headers (bold,
# centered, bordered).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New header labels.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Team record is constant for every player row (2-45): 98 wins, 64 losses,
# 0 ties.
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 29).Value = 98   # AC
    $ws.Cells.Item($row, 30).Value = 64   # AD
    $ws.Cells.Item($row, 31).Value = 0    # AE
}
